$wb = $excel.ActiveWorkbook

$wsDaily = $wb.Worksheets.Item("Tagesergebnisse")
$wsMonthly = $wb.Worksheets.Item("Monatsergebnisse")
$wsTotal = $wb.Worksheets.Item("Gesamtergebnis")

# Gesamtergebnis ("Total") sheet: fill in totals for the balance columns (C3/D3)
# which were previously "N/A" text, now numeric totals matching row 2.
$wsTotal.Range("C3").Value = 1.19
$wsTotal.Range("D3").Value = 0.77

# Update selections to match the saved workbook state (active cell moved to C4).
[void]$wsTotal.Range("C4").Select()

[void]$wsDaily.Range("I1").Select()
[void]$wsMonthly.Range("H10").Select()

[void]$wsTotal.Activate()
